# Updated cryptos list on Sun Oct 15 21:37:21 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, $Text)
    # Force the cell to be stored as text (avoids Excel auto-converting
    # number-like strings such as "0.493" or "211.09" into numeric values),
    # then restore the default "Normal" style so no stray formatting is left
    # behind on the cell.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell $ws.Range("D2") "27.335.11"

# Row 3 - Ethereum
Set-TextCell $ws.Range("D3") "1.569.24"
$ws.Range("E3").Value = "  +0.46%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
Set-TextCell $ws.Range("D5") "211.09"
$ws.Range("E5").Value = "  +1.72%  "

# Row 6 - XRP
Set-TextCell $ws.Range("D6") "0.493"
$ws.Range("E6").Value = "  +0.64%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - Solana
Set-TextCell $ws.Range("D8") "22.16"
$ws.Range("E8").Value = "  +0.44%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.65%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.10%  "

# Row 11 - TRON
Set-TextCell $ws.Range("D11") "0.0870"
$ws.Range("E11").Value = "  +1.54%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell $ws.Range("D12") "1.792.17"
$ws.Range("E12").Value = "  +0.42%  "

# Row 13 - WrappedEther
Set-TextCell $ws.Range("D13") "1.564.25"
$ws.Range("E13").Value = "  -0.02%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.80%  "

# Row 15 - Polygon
Set-TextCell $ws.Range("D15") "0.519"
$ws.Range("E15").Value = "  +0.00%  "

# Row 16 - WrappedBTC
Set-TextCell $ws.Range("D16") "27.271.66"

# Row 17 - Litecoin
Set-TextCell $ws.Range("D17") "62.27"
$ws.Range("E17").Value = "  +0.30%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  +2.31%  "

# Row 19 - BitcoinCash
Set-TextCell $ws.Range("D19") "218.03"
$ws.Range("E19").Value = "  +0.70%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -0.21%  "

# Row 22 - Uniswap
Set-TextCell $ws.Range("D22") "4.16"
$ws.Range("E22").Value = "  +1.25%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +0.28%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.09%  "

# Row 25 - Monero
Set-TextCell $ws.Range("D25") "153.62"
$ws.Range("E25").Value = "  +0.69%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +0.67%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +0.19%  "

# Row 28 - Stellar
Set-TextCell $ws.Range("D28") "0.107"
$ws.Range("E28").Value = "  +1.89%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.03%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +2.74%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +0.28%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.51%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +1.84%  "

# Row 34 - Maker
Set-TextCell $ws.Range("D34") "1.446.72"
$ws.Range("E34").Value = "  +1.82%  "

# Row 35 - TrustWalletToken
Set-TextCell $ws.Range("D35") "1.10"
$ws.Range("E35").Value = "  +2.57%  "

# Row 36 - LidoDAOToken
Set-TextCell $ws.Range("D36") "1.61"

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +0.35%  "

# Row 38 - VeChain
Set-TextCell $ws.Range("D38") "0.0167"
$ws.Range("E38").Value = "  +0.95%  "

# Row 39 - ImmutableX
Set-TextCell $ws.Range("D39") "0.534"
$ws.Range("E39").Value = "  +0.16%  "

# Row 40 - FraxShare
Set-TextCell $ws.Range("D40") "5.91"
$ws.Range("E40").Value = "  +2.15%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +0.33%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.09%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +0.58%  "

# Row 44 - WEMIXToken
$ws.Range("E44").Value = "  -0.43%  "

# Row 45 - Aave
Set-TextCell $ws.Range("D45") "64.67"
$ws.Range("E45").Value = "  -0.25%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  -0.70%  "

# Row 47 - RocketPoolETH
Set-TextCell $ws.Range("D47") "1.704.17"
$ws.Range("E47").Value = "  +0.31%  "

# Row 48 - Quant
Set-TextCell $ws.Range("D48") "86.09"
$ws.Range("E48").Value = "  -1.48%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +1.10%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  +0.17%  "

# Row 51 - Algorand
Set-TextCell $ws.Range("D51") "0.0956"
$ws.Range("E51").Value = "  -0.32%  "
